$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 02:20"

# Chequia overtakes Polonia in the ranking -> swap country labels at rows 32/33
$ws.Range("A32").Value = "Chequia"
$ws.Range("A33").Value = "Polonia"

# Refreshed COVID stats per country (Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes)
# Row 4
$ws.Range("B4").Value = 8519467
$ws.Range("C4").Value = 60589
$ws.Range("D4").Value = 5543053
$ws.Range("E4").Value = 2750294
$ws.Range("G4").Value = 888
$ws.Range("H4").Value = 226120

# Row 5
$ws.Range("B5").Value = 7649158
$ws.Range("C5").Value = 54422
$ws.Range("D5").Value = 6792550
$ws.Range("E5").Value = 740658
$ws.Range("G5").Value = 714
$ws.Range("H5").Value = 115950

# Row 6
$ws.Range("B6").Value = 5274817
$ws.Range("C6").Value = 23690
$ws.Range("E6").Value = 398336
$ws.Range("G6").Value = 662
$ws.Range("H6").Value = 154888

# Row 9
$ws.Range("B9").Value = 1018999
$ws.Range("C9").Value = 16337
$ws.Range("D9").Value = 829647
$ws.Range("E9").Value = 162252
$ws.Range("G9").Value = 384
$ws.Range("H9").Value = 27100

# Row 21
$ws.Range("B21").Value = 380898
$ws.Range("C21").Value = 7167
$ws.Range("E21").Value = 72643

# Row 31
$ws.Range("B31").Value = 203688
$ws.Range("C31").Value = 2251
$ws.Range("D31").Value = 171743
$ws.Range("E31").Value = 22151
$ws.Range("G31").Value = 16
$ws.Range("H31").Value = 9794

# Row 32
$ws.Range("B32").Value = 193946
$ws.Range("C32").Value = 11984
$ws.Range("D32").Value = 79108
$ws.Range("E32").Value = 113219
$ws.Range("G32").Value = 106
$ws.Range("H32").Value = 1619

# Row 33
$ws.Range("B33").Value = 192539
$ws.Range("C33").Value = 9291
$ws.Range("D33").Value = 95956
$ws.Range("E33").Value = 92862
$ws.Range("G33").Value = 107
$ws.Range("H33").Value = 3721

# Row 40
$ws.Range("B40").Value = 125739
$ws.Range("C40").Value = 558
$ws.Range("D40").Value = 102028
$ws.Range("E40").Value = 21126
$ws.Range("G40").Value = 11
$ws.Range("H40").Value = 2585

# Row 55
$ws.Range("B55").Value = 87644
$ws.Range("C55").Value = 483
$ws.Range("E55").Value = 6581
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 747

# Row 66
$ws.Range("B66").Value = 56073
$ws.Range("C66").Value = 621
$ws.Range("D66").Value = 37167
$ws.Range("E66").Value = 17675
$ws.Range("G66").Value = 24
$ws.Range("H66").Value = 1231

# Row 121
$ws.Range("B121").Value = 7329
$ws.Range("C121").Value = 207
$ws.Range("E121").Value = 5015
$ws.Range("G121").Value = 19
$ws.Range("H121").Value = 115

# Row 126
$ws.Range("B126").Value = 5800
$ws.Range("C126").Value = 12
$ws.Range("D126").Value = 5437
$ws.Range("E126").Value = 247

# Row 134
$ws.Range("B134").Value = 5144
$ws.Range("C134").Value = 11
$ws.Range("D134").Value = 4979
$ws.Range("E134").Value = 56

# Row 155
$ws.Range("B155").Value = 2623
$ws.Range("C155").Value = 63
$ws.Range("D155").Value = 2142
$ws.Range("E155").Value = 429
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 52
